$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 770
$ws.Range("I2").Value = 770
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 770
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -657
$ws.Range("N2").ClearContents()
$ws.Range("H9").Value = 1681.9
$ws.Range("I9").Value = 175.5
$ws.Range("J9").Value = 7707.5
$ws.Range("K9").Value = 175.5
$ws.Range("L9").Value = 7707.5
$ws.Range("M9").Value = -6.5
$ws.Range("N9").Value = -8045.5
$ws.Range("H12").Value = 316
$ws.Range("I12").Value = 316
$ws.Range("K12").Value = 316
$ws.Range("M12").Value = -146
$ws.Range("H21").Value = 9650.272000000001
$ws.Range("I21").Value = 9650.272000000001
$ws.Range("K21").Value = 9650.272000000001
$ws.Range("M21").Value = -9182.272000000001
$ws.Range("H23").Value = 9650.272000000001
$ws.Range("I23").Value = 9650.272000000001
$ws.Range("K23").Value = 9650.272000000001
$ws.Range("M23").Value = -9416.272000000001
$ws.Range("H32").Value = 4636.5557
$ws.Range("I32").Value = 506.66666
$ws.Range("K32").Value = 506.66666
$ws.Range("M32").Value = -180.66666
$ws.Range("H38").Value = 541.1539
$ws.Range("I38").Value = 541.1539
$ws.Range("K38").Value = 1623.4617
$ws.Range("M38").Value = -1251.4617
$ws.Range("H49").Value = 2160
$ws.Range("I49").Value = 2160
$ws.Range("K49").Value = 6480
$ws.Range("M49").Value = -6344
$ws.Range("H58").Value = 1817.75
$ws.Range("J58").Value = 4012.75
$ws.Range("L58").Value = 12038.25
$ws.Range("N58").Value = -12338.25
$ws.Range("H62").Value = 3999.2
$ws.Range("I62").Value = 3999.1428
$ws.Range("K62").Value = 3999.1428
$ws.Range("M62").Value = -3375.1428
$ws.Range("H64").Value = 5172.276
$ws.Range("I64").Value = 5192.3
$ws.Range("J64").Value = 5127.778
$ws.Range("K64").Value = 5192.3
$ws.Range("L64").Value = 5127.778
$ws.Range("M64").Value = -4944.3
$ws.Range("N64").Value = -5623.778
$ws.Range("H65").Value = 3999.2
$ws.Range("I65").Value = 3999.1428
$ws.Range("K65").Value = 19995.714
$ws.Range("M65").Value = -16875.714
$ws.Range("H67").Value = 5172.276
$ws.Range("I67").Value = 5192.3
$ws.Range("J67").Value = 5127.778
$ws.Range("K67").Value = 5192.3
$ws.Range("L67").Value = 5127.778
$ws.Range("M67").Value = -4334.3
$ws.Range("N67").Value = -6843.778
$ws.Range("H70").Value = 3088.25
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3088.25
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 9264.75
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -9804.75
$ws.Range("H73").Value = 3088.25
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3088.25
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 9264.75
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -11136.75
$ws.Range("H74").Value = 6060.4375
$ws.Range("I74").Value = 5711.643
$ws.Range("K74").Value = 5711.643
$ws.Range("M74").Value = -4775.643
$ws.Range("H77").Value = 6060.4375
$ws.Range("I77").Value = 5711.643
$ws.Range("K77").Value = 28558.215
$ws.Range("M77").Value = -23878.215
$ws.Range("H98").Value = 1322.1111
$ws.Range("I98").Value = 1206.1875
$ws.Range("J98").Value = 2249.5
$ws.Range("K98").Value = 1206.1875
$ws.Range("L98").Value = 2249.5
$ws.Range("M98").Value = 291.8125
$ws.Range("N98").Value = -5245.5
$ws.Range("H99").Value = 2694.5
$ws.Range("I99").Value = 1708.5
$ws.Range("K99").Value = 5125.5
$ws.Range("M99").Value = -3627.5
$ws.Range("H106").Value = 1331
$ws.Range("I106").Value = 1214.1
$ws.Range("K106").Value = 1214.1
$ws.Range("M106").Value = -583.0999999999999
$ws.Range("H115").Value = 1471.5
$ws.Range("J115").Value = 1000
$ws.Range("L115").Value = 3000
$ws.Range("N115").Value = -6134
$ws.Range("H122").Value = 1322.1111
$ws.Range("I122").Value = 1206.1875
$ws.Range("J122").Value = 2249.5
$ws.Range("K122").Value = 3618.5625
$ws.Range("L122").Value = 6748.5
$ws.Range("M122").Value = -1168.5625
$ws.Range("N122").Value = -11648.5
$ws.Range("H132").Value = 65105.625
$ws.Range("I132").Value = 72985.5
$ws.Range("K132").Value = 218956.5
$ws.Range("M132").Value = -216426.5
$ws.Range("H138").Value = 2357.44
$ws.Range("I138").Value = 1630.7333
$ws.Range("J138").Value = 3447.5
$ws.Range("K138").Value = 4892.199900000001
$ws.Range("L138").Value = 10342.5
$ws.Range("M138").Value = 247.8000999999995
$ws.Range("N138").Value = -20622.5
$ws.Range("H141").Value = 3597.3076
$ws.Range("I141").Value = 3597.3076
$ws.Range("K141").Value = 10791.9228
$ws.Range("M141").Value = -5611.9228

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2835
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 4250
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 4250
$ws.Range("M8").Value = 139
$ws.Range("N8").Value = -4538
$ws.Range("H37").Value = 2166.6667
$ws.Range("J37").Value = 2500
$ws.Range("L37").Value = 2500
$ws.Range("N37").Value = -3046
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H61").Value = 3772.9333
$ws.Range("I61").Value = 3177.111
$ws.Range("J61").Value = 4666.6665
$ws.Range("K61").Value = 3177.111
$ws.Range("L61").Value = 4666.6665
$ws.Range("M61").Value = -2965.111
$ws.Range("N61").Value = -5090.6665
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 622
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 3628.111
$ws.Range("I122").Value = 4379.5
$ws.Range("J122").Value = 2125.3333
$ws.Range("K122").Value = 13138.5
$ws.Range("L122").Value = 6375.999899999999
$ws.Range("M122").Value = -10688.5
$ws.Range("N122").Value = -11275.9999
$ws.Range("H132").Value = 35716780
$ws.Range("I132").Value = 2617.1667
$ws.Range("J132").Value = 250001740
$ws.Range("K132").Value = 7851.500100000001
$ws.Range("L132").Value = 750005220
$ws.Range("M132").Value = -5321.500100000001
$ws.Range("N132").Value = -750010280
$ws.Range("H136").Value = 3772.9333
$ws.Range("I136").Value = 3177.111
$ws.Range("J136").Value = 4666.6665
$ws.Range("K136").Value = 9531.332999999999
$ws.Range("L136").Value = 13999.9995
$ws.Range("M136").Value = -6981.332999999999
$ws.Range("N136").Value = -19099.9995
$ws.Range("H139").Value = 35715
$ws.Range("J139").Value = 35715
$ws.Range("L139").Value = 35715
$ws.Range("N139").Value = -45995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H10").Value = 3998
$ws.Range("I10").Value = 3998
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 3998
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -3858
$ws.Range("N10").ClearContents()
$ws.Range("H14").Value = 5265.3335
$ws.Range("J14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("N14").Value = -15344
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1711
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H86").Value = 3336.1428
$ws.Range("I86").Value = 2586.75
$ws.Range("K86").Value = 2586.75
$ws.Range("M86").Value = -1463.75
$ws.Range("H89").Value = 3336.1428
$ws.Range("I89").Value = 2586.75
$ws.Range("K89").Value = 12933.75
$ws.Range("M89").Value = -7317.75
$ws.Range("H132").Value = 69816.664
$ws.Range("J132").Value = 69816.664
$ws.Range("L132").Value = 69816.664
$ws.Range("N132").Value = -79936.664
$ws.Range("H134").Value = 44880480
$ws.Range("I134").Value = 22737172
$ws.Range("J134").Value = 166668670
$ws.Range("K134").Value = 68211516
$ws.Range("L134").Value = 500006010
$ws.Range("M134").Value = -68208981
$ws.Range("N134").Value = -500011080
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2244
$ws.Range("J16").Value = 2918.2
$ws.Range("L16").Value = 2918.2
$ws.Range("N16").Value = -3492.2
$ws.Range("H31").Value = 2504.875
$ws.Range("I31").Value = 2097.25
$ws.Range("K31").Value = 2097.25
$ws.Range("M31").Value = -1802.25
$ws.Range("H34").Value = 2504.875
$ws.Range("I34").Value = 2097.25
$ws.Range("K34").Value = 2097.25
$ws.Range("M34").Value = -1895.25
$ws.Range("H58").Value = 1881.8235
$ws.Range("I58").Value = 1612.7333
$ws.Range("J58").Value = 3900
$ws.Range("K58").Value = 1612.7333
$ws.Range("L58").Value = 3900
$ws.Range("M58").Value = -1409.7333
$ws.Range("N58").Value = -4306
$ws.Range("H62").Value = 10166.706
$ws.Range("I62").Value = 3757.4614
$ws.Range("K62").Value = 3757.4614
$ws.Range("M62").Value = -3133.4614
$ws.Range("H65").Value = 10166.706
$ws.Range("I65").Value = 3757.4614
$ws.Range("K65").Value = 18787.307
$ws.Range("M65").Value = -15667.307
$ws.Range("H99").Value = 2245.6365
$ws.Range("I99").Value = 2057.8333
$ws.Range("K99").Value = 2057.8333
$ws.Range("M99").Value = -559.8332999999998
$ws.Range("H105").Value = 2547.3044
$ws.Range("I105").Value = 2047.3125
$ws.Range("K105").Value = 2047.3125
$ws.Range("M105").Value = -300.3125
$ws.Range("H113").Value = 2244
$ws.Range("J113").Value = 2918.2
$ws.Range("L113").Value = 2918.2
$ws.Range("N113").Value = -7258.2
$ws.Range("H122").Value = 1693.1052
$ws.Range("I122").Value = 935.2857
$ws.Range("J122").Value = 2135.1667
$ws.Range("K122").Value = 2805.8571
$ws.Range("L122").Value = 6405.500100000001
$ws.Range("M122").Value = -355.8571000000002
$ws.Range("N122").Value = -11305.5001
$ws.Range("H126").Value = 2245.6365
$ws.Range("I126").Value = 2057.8333
$ws.Range("K126").Value = 6173.499899999999
$ws.Range("M126").Value = -3703.499899999999
$ws.Range("H132").Value = 4291.073
$ws.Range("I132").Value = 4180.641
$ws.Range("K132").Value = 12541.923
$ws.Range("M132").Value = -10011.923
$ws.Range("H136").Value = 1881.8235
$ws.Range("I136").Value = 1612.7333
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 4838.199900000001
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -2288.199900000001
$ws.Range("N136").Value = -16800
$ws.Range("H139").Value = 40709
$ws.Range("I139").Value = 40709
$ws.Range("K139").Value = 40709
$ws.Range("M139").Value = -35569

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 377.6842
$ws.Range("I7").Value = 388.54544
$ws.Range("J7").Value = 362.75
$ws.Range("K7").Value = 1165.63632
$ws.Range("L7").Value = 1088.25
$ws.Range("M7").Value = -1053.63632
$ws.Range("N7").Value = -1312.25
$ws.Range("H13").Value = 1405.5
$ws.Range("I13").Value = 2313.5
$ws.Range("J13").Value = 497.5
$ws.Range("K13").Value = 6940.5
$ws.Range("L13").Value = 1492.5
$ws.Range("M13").Value = -6772.5
$ws.Range("N13").Value = -1828.5
$ws.Range("H23").Value = 913.1905
$ws.Range("I23").Value = 1040
$ws.Range("K23").Value = 3120
$ws.Range("M23").Value = -2885
$ws.Range("H34").Value = 5164.7646
$ws.Range("J34").Value = 5574.7334
$ws.Range("L34").Value = 16724.2002
$ws.Range("N34").Value = -16892.2002
$ws.Range("H37").Value = 120000
$ws.Range("J37").Value = 120000
$ws.Range("L37").Value = 360000
$ws.Range("N37").Value = -360224
$ws.Range("H39").Value = 4542.0527
$ws.Range("J39").Value = 5153.467
$ws.Range("L39").Value = 15460.401
$ws.Range("N39").Value = -16048.401
$ws.Range("H49").Value = 3833.3
$ws.Range("I49").Value = 3624.9583
$ws.Range("J49").Value = 4666.6665
$ws.Range("K49").Value = 10874.8749
$ws.Range("L49").Value = 13999.9995
$ws.Range("M49").Value = -10718.8749
$ws.Range("N49").Value = -14311.9995
$ws.Range("H55").Value = 2471.8
$ws.Range("J55").Value = 2988.75
$ws.Range("L55").Value = 8966.25
$ws.Range("N55").Value = -9320.25
$ws.Range("H63").Value = 1009
$ws.Range("I63").Value = 1009
$ws.Range("K63").Value = 3027
$ws.Range("M63").Value = -2278
$ws.Range("H64").Value = 13771
$ws.Range("J64").Value = 20007
$ws.Range("L64").Value = 60021
$ws.Range("N64").Value = -60561
$ws.Range("H66").Value = 1009
$ws.Range("I66").Value = 1009
$ws.Range("K66").Value = 9081
$ws.Range("M66").Value = -5337
$ws.Range("H67").Value = 13771
$ws.Range("J67").Value = 20007
$ws.Range("L67").Value = 60021
$ws.Range("N67").Value = -61893
$ws.Range("H75").Value = 10580
$ws.Range("J75").Value = 10975
$ws.Range("L75").Value = 32925
$ws.Range("N75").Value = -34921
$ws.Range("H78").Value = 10580
$ws.Range("J78").Value = 10975
$ws.Range("L78").Value = 98775
$ws.Range("N78").Value = -108759
$ws.Range("H87").Value = 27500
$ws.Range("I87").Value = 20000
$ws.Range("K87").Value = 60000
$ws.Range("M87").Value = -58752
$ws.Range("H90").Value = 27500
$ws.Range("I90").Value = 20000
$ws.Range("K90").Value = 180000
$ws.Range("M90").Value = -173760
$ws.Range("H95").Value = 34162.668
$ws.Range("I95").Value = 24988
$ws.Range("K95").Value = 74964
$ws.Range("M95").Value = -72905
$ws.Range("H103").Value = 211.38461
$ws.Range("I103").Value = 230.28572
$ws.Range("K103").Value = 690.85716
$ws.Range("M103").Value = 188.14284
$ws.Range("H113").Value = 798.625
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 848.1667
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 2544.5001
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -6884.5001
$ws.Range("H116").Value = 81029.164
$ws.Range("I116").Value = 90095.31
$ws.Range("K116").Value = 270285.93
$ws.Range("M116").Value = -266843.93
$ws.Range("H117").Value = 3009.6924
$ws.Range("I117").Value = 900
$ws.Range("J117").Value = 3185.5
$ws.Range("K117").Value = 2700
$ws.Range("L117").Value = 9556.5
$ws.Range("M117").Value = 742
$ws.Range("N117").Value = -16440.5
$ws.Range("H129").Value = 2583.8572
$ws.Range("I129").Value = 2196
$ws.Range("J129").Value = 2874.75
$ws.Range("K129").Value = 6588
$ws.Range("L129").Value = 8624.25
$ws.Range("M129").Value = -1588
$ws.Range("N129").Value = -18624.25
$ws.Range("H131").Value = 569554
$ws.Range("I131").Value = 1041.2142
$ws.Range("K131").Value = 3123.6426
$ws.Range("M131").Value = 1916.3574
$ws.Range("H140").Value = 1281.2174
$ws.Range("I140").Value = 1203.091
$ws.Range("K140").Value = 3609.273
$ws.Range("M140").Value = 1570.727
$ws.Range("H141").Value = 1899
$ws.Range("I141").Value = 1899
$ws.Range("K141").Value = 5697
$ws.Range("M141").Value = -517

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H43").Value = 670633.3
$ws.Range("I43").Value = 670633.3
$ws.Range("K43").Value = 670633.3
$ws.Range("M43").Value = -670482.3
$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5312
$ws.Range("H62").Value = 25042.25
$ws.Range("J62").Value = 25042.25
$ws.Range("L62").Value = 25042.25
$ws.Range("N62").Value = -26414.25
$ws.Range("H65").Value = 25042.25
$ws.Range("J65").Value = 25042.25
$ws.Range("L65").Value = 75126.75
$ws.Range("N65").Value = -81990.75
$ws.Range("H70").Value = 11144.77
$ws.Range("I70").Value = 11088.7
$ws.Range("J70").Value = 11331.667
$ws.Range("K70").Value = 11088.7
$ws.Range("L70").Value = 11331.667
$ws.Range("M70").Value = -10818.7
$ws.Range("N70").Value = -11871.667
$ws.Range("H73").Value = 11144.77
$ws.Range("I73").Value = 11088.7
$ws.Range("J73").Value = 11331.667
$ws.Range("K73").Value = 11088.7
$ws.Range("L73").Value = 11331.667
$ws.Range("M73").Value = -10152.7
$ws.Range("N73").Value = -13203.667
$ws.Range("H122").Value = 4150.2144
$ws.Range("I122").Value = 4827.737
$ws.Range("J122").Value = 2719.889
$ws.Range("K122").Value = 14483.211
$ws.Range("L122").Value = 8159.667
$ws.Range("M122").Value = -12033.211
$ws.Range("N122").Value = -13059.667
$ws.Range("H136").Value = 39899.5
$ws.Range("J136").Value = 39899.5
$ws.Range("L136").Value = 119698.5
$ws.Range("N136").Value = -124798.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1615
$ws.Range("I7").Value = 1615
$ws.Range("K7").Value = 1615
$ws.Range("M7").Value = -1503
$ws.Range("H40").Value = 4159
$ws.Range("I40").Value = 4176.6665
$ws.Range("K40").Value = 4176.6665
$ws.Range("M40").Value = -4040.6665
$ws.Range("H55").Value = 285.73334
$ws.Range("I55").Value = 311.75
$ws.Range("K55").Value = 311.75
$ws.Range("M55").Value = -138.75
$ws.Range("H64").Value = 49997.57
$ws.Range("J64").Value = 49997.168
$ws.Range("L64").Value = 49997.168
$ws.Range("N64").Value = -50447.168
$ws.Range("H67").Value = 49997.57
$ws.Range("J67").Value = 49997.168
$ws.Range("L67").Value = 49997.168
$ws.Range("N67").Value = -51557.168
$ws.Range("H68").Value = 3528.25
$ws.Range("I68").Value = 3421.077
$ws.Range("J68").Value = 3992.6667
$ws.Range("K68").Value = 3421.077
$ws.Range("L68").Value = 3992.6667
$ws.Range("M68").Value = -2672.077
$ws.Range("N68").Value = -5490.6667
$ws.Range("H71").Value = 3528.25
$ws.Range("I71").Value = 3421.077
$ws.Range("J71").Value = 3992.6667
$ws.Range("K71").Value = 17105.385
$ws.Range("L71").Value = 19963.3335
$ws.Range("M71").Value = -13361.385
$ws.Range("N71").Value = -27451.3335
$ws.Range("H122").Value = 3492.7273
$ws.Range("I122").Value = 3330
$ws.Range("J122").Value = 3528.889
$ws.Range("K122").Value = 9990
$ws.Range("L122").Value = 10586.667
$ws.Range("M122").Value = -7540
$ws.Range("N122").Value = -15486.667
$ws.Range("H126").Value = 1615
$ws.Range("I126").Value = 1615
$ws.Range("K126").Value = 4845
$ws.Range("M126").Value = -2375
$ws.Range("H132").Value = 9915.429
$ws.Range("I132").Value = 9901.5
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 29704.5
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -27174.5
$ws.Range("N132").Value = -35057

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 25029500
$ws.Range("J5").Value = 25029500
$ws.Range("L5").Value = 25029500
$ws.Range("N5").Value = -25029724
$ws.Range("H54").Value = 18848.268
$ws.Range("J54").Value = 18848.268
$ws.Range("L54").Value = 18848.268
$ws.Range("N54").Value = -19888.268
$ws.Range("H100").Value = 1812.5454
$ws.Range("I100").Value = 1548.3334
$ws.Range("K100").Value = 3096.6668
$ws.Range("M100").Value = -2555.6668
$ws.Range("H126").Value = 2258.0715
$ws.Range("I126").Value = 1625
$ws.Range("J126").Value = 2891.1428
$ws.Range("K126").Value = 4875
$ws.Range("L126").Value = 8673.428400000001
$ws.Range("M126").Value = -2405
$ws.Range("N126").Value = -13613.4284
$ws.Range("H132").Value = 1621.6842
$ws.Range("I132").Value = 1548.1765
$ws.Range("K132").Value = 4644.529500000001
$ws.Range("M132").Value = -2114.529500000001
$ws.Range("H136").Value = 1899.75
$ws.Range("I136").Value = 1712.375
$ws.Range("J136").Value = 2274.5
$ws.Range("K136").Value = 5137.125
$ws.Range("L136").Value = 6823.5
$ws.Range("M136").Value = -2587.125
$ws.Range("N136").Value = -11923.5

